$wb = $excel.ActiveWorkbook

function Set-HeaderStyle($rng) {
    $rng.Font.Bold = $true
    $rng.HorizontalAlignment = -4108   # xlCenter
    $rng.VerticalAlignment = -4160     # xlTop
    $rng.Borders.LineStyle = 1
}

# --- New sheet: area_lores_basic (appended after last sheet) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "area_lores_basic"

$ws3.Range("A1").Value = "index"
$ws3.Range("B1").Value = "area"
Set-HeaderStyle $ws3.Range("A1:B1")

$ws3.Range("A2").Value = "count"
$ws3.Range("B2").Value = 14

$ws3.Range("A3").Value = "mean"
$ws3.Range("B3").Value = 14.40594373249588

$ws3.Range("A4").Value = "std"
$ws3.Range("B4").Value = 15.43606997512803

$ws3.Range("A5").Value = "min"
$ws3.Range("B5").Value = 2.54798191151668

$ws3.Range("A6").NumberFormat = "@"
$ws3.Range("A6").Value = "25%"
$ws3.Range("A6").Style = "Normal"
$ws3.Range("B6").Value = 4.657497179950208

$ws3.Range("A7").NumberFormat = "@"
$ws3.Range("A7").Value = "50%"
$ws3.Range("A7").Style = "Normal"
$ws3.Range("B7").Value = 9.097167038340388

$ws3.Range("A8").NumberFormat = "@"
$ws3.Range("A8").Value = "75%"
$ws3.Range("A8").Style = "Normal"
$ws3.Range("B8").Value = 18.15870380406473

$ws3.Range("A9").Value = "max"
$ws3.Range("B9").Value = 60.53123548751798

# --- New sheet: area_pop_sum_basic (appended after area_lores_basic) ---
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "area_pop_sum_basic"

$ws4.Range("A1").Value = "index"
$ws4.Range("B1").Value = 0
Set-HeaderStyle $ws4.Range("A1:B1")

$ws4.Range("A2").Value = "area"
$ws4.Range("B2").Value = 201.6832122549423

$ws4.Range("A3").Value = "population"
$ws4.Range("B3").Value = 228201

$ws4.Range("A4").Value = "density"
$ws4.Range("B4").Value = 1131.482375000738

# Restore original active sheet (first sheet) so the new sheets don't
# steal the workbook's active-tab focus.
$wb.Worksheets.Item(1).Activate()
